# Add a new "canonical SMILES" column (D) to the microstate list sheet.
# For microstates whose "canonical isomeric SMILES" (column C) encodes E/Z
# bond stereochemistry with '/' and '\' markers, the new column gives the
# same SMILES with those directional-bond markers stripped. For microstates
# that have no such markers, column D simply repeats column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("D2").Value2 = "canonical SMILES"

# Data rows
$ws.Range("D3").Value2  = "c1cc(cc(c1)Br)Nc2[nH]cc(c(=[NH+]c3cccc(c3)Br)n2)F"
$ws.Range("D4").Value2  = $ws.Range("C4").Value2
$ws.Range("D5").Value2  = "c1cc(cc(c1)Br)Nc2c(cnc(=Nc3cccc(c3)Br)[nH]2)F"
$ws.Range("D6").Value2  = $ws.Range("C6").Value2
$ws.Range("D7").Value2  = "c1cc(cc(c1)Br)Nc2[nH]c(=Nc3cccc(c3)Br)c(cn2)F"
$ws.Range("D8").Value2  = $ws.Range("C8").Value2
$ws.Range("D9").Value2  = "c1cc(cc(c1)Br)Nc2[nH]cc(c(=Nc3cccc(c3)Br)n2)F"
$ws.Range("D10").Value2 = "c1cc(cc(c1)Br)[N-]c2[nH]c(=Nc3cccc(c3)Br)c(cn2)F"
$ws.Range("D11").Value2 = "c1cc(cc(c1)Br)Nc2c(c[nH]c(=Nc3cccc(c3)Br)n2)F"
$ws.Range("D12").Value2 = $ws.Range("C12").Value2
$ws.Range("D13").Value2 = $ws.Range("C13").Value2
$ws.Range("D14").Value2 = $ws.Range("C14").Value2
$ws.Range("D15").Value2 = "c1cc(cc(c1)Br)[N-]c2[nH]cc(c(=Nc3cccc(c3)Br)n2)F"
$ws.Range("D16").Value2 = "c1cc(cc(c1)Br)Nc2c(c[nH]c(=[NH+]c3cccc(c3)Br)[nH+]2)F"
$ws.Range("D17").Value2 = "c1cc(cc(c1)Br)[NH2+]c2c(c[nH]c(=[NH+]c3cccc(c3)Br)n2)F"
$ws.Range("D18").Value2 = $ws.Range("C18").Value2
$ws.Range("D19").Value2 = $ws.Range("C19").Value2
$ws.Range("D20").Value2 = $ws.Range("C20").Value2
$ws.Range("D21").Value2 = $ws.Range("C21").Value2
$ws.Range("D22").Value2 = $ws.Range("C22").Value2
$ws.Range("D23").Value2 = $ws.Range("C23").Value2
$ws.Range("D24").Value2 = $ws.Range("C24").Value2
$ws.Range("D25").Value2 = $ws.Range("C25").Value2
$ws.Range("D26").Value2 = $ws.Range("C26").Value2

# Give the new column a sensible custom width, matching the other
# SMILES/description columns on this sheet.
$ws.Range("D1").ColumnWidth = 43.3
